# Applies two changes to the presentation:
#   1. Re-styles the single table (slide 16) to the built-in
#      "{312F6497-F09B-4C69-BFC9-78C70A7F5908}" table style.
#   2. Switches the deck's theme colour scheme from the custom
#      "Integral" palette to the default "Office Theme" palette
#      (dk2/lt2/accent1-6/hlink/folHlink - dk1/lt1 are already
#      identical black/white in both palettes).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$newTableStyleId = "{312F6497-F09B-4C69-BFC9-78C70A7F5908}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Theme colours -------------------------------------------------
# Colors() index order exposed by the theme colour scheme object:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $rgb = $officeThemeColors[$i - 1]
    $r = $rgb[0]
    $g = $rgb[1]
    $b = $rgb[2]
    $colorScheme.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
